# Inserts a new market-day block (3 rows: Maduro / Pintón / Primera Pintón,
# dated 2021-11-04) right after the existing row 440, pushing all the
# subsequent rows down by 3 (old row 441 -> new row 444, ..., old row 513 ->
# new row 516). The dimension of the sheet grows from A1:T513 to A1:T516
# automatically as a consequence of the row insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert 3 blank rows before row 441 (shifts 441:513 -> 444:516) ---
$ws.Rows("441:443").Insert()

# --- 2. Populate the 3 newly inserted rows with the new observations ---

$rows = @(
    @{ Row = 441; Date = 44504; Quality = "Maduro";         Volume = 160; Min = 15000; Max = 15000; Avg = 15000; PerKg = 750 },
    @{ Row = 442; Date = 44504; Quality = "Pintón";         Volume = 320; Min = 17000; Max = 17000; Avg = 17000; PerKg = 850 },
    @{ Row = 443; Date = 44504; Quality = "Primera Pintón"; Volume = 400; Min = 18000; Max = 19000; Avg = 18600; PerKg = 930 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value  = 3                                   # A Mercado ID
    $ws.Cells.Item($row, 2).Value  = "Femacal de La Calera"               # B Mercado
    $ws.Cells.Item($row, 3).Value  = "Coquimbo"                           # C Región
    $ws.Cells.Item($row, 4).Value  = $r.Date                              # D Fecha
    $ws.Cells.Item($row, 5).Value  = 5                                   # E Codreg
    $ws.Cells.Item($row, 6).Value  = "Fruta"                              # F Tipo
    $ws.Cells.Item($row, 7).Value  = 100108                               # G Producto ID
    $ws.Cells.Item($row, 8).Value  = "Tropicales y subtropicales"         # H Producto
    $ws.Cells.Item($row, 9).Value  = 100108006                            # I Categoría ID
    $ws.Cells.Item($row, 10).Value = "Plátano"                            # J Categoría
    $ws.Cells.Item($row, 11).Value = "Sin especificar"                   # K Variedad
    $ws.Cells.Item($row, 12).Value = $r.Quality                          # L Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volume                           # M Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min                              # N Precio mínimo
    $ws.Cells.Item($row, 15).Value = $r.Max                              # O Precio máximo
    $ws.Cells.Item($row, 16).Value = $r.Avg                              # P Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = "$/caja 20 kilos"                   # Q Unidad de comercialización
    $ws.Cells.Item($row, 18).Value = "Ecuador"                           # R Origen
    $ws.Cells.Item($row, 19).Value = $r.PerKg                            # S Precio $/Kg
    $ws.Cells.Item($row, 20).Value = 20                                  # T Kg / unidad
}
